# Generate Report for Handoff
# Updates the localization-status workbook with a newly generated handoff
# package: new guid-named source file, new xliff hashes, and refreshed
# timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "159ab9b9-fee5-42bb-9b23-80153ebdbf90"
$newGuid = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb"

$oldZhXlf = "159ab9b9-fee5-42bb-9b23-80153ebdbf90.136f23df5424049d121cf2498532e31cf47e2356.zh-cn.xlf"
$newZhXlf = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.6cde7396a0c1c32beacae51b5b1dd5485bd2aa91.zh-cn.xlf"

$oldDeXlf = "159ab9b9-fee5-42bb-9b23-80153ebdbf90.136f23df5424049d121cf2498532e31cf47e2356.de-de.xlf"
$newDeXlf = "849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.6cde7396a0c1c32beacae51b5b1dd5485bd2aa91.de-de.xlf"

$oldGenDate = "2016-09-05 07:07:15"
$newGenDate = "2016-09-05 07:07:32"

$oldZhDate = "2016-09-05 07:07:09"
$newZhDate = "2016-09-05 07:07:27"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newGenDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhDate

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenDate
